$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 0.625
$ws.Range("F2").Value = 0.6959323534415676
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0.9999999999

# Row 3
$ws.Range("E3").Value = 0.4858680360343469
$ws.Range("G3").Value = 1
$ws.Range("I3").Value = 1

# Row 4
$ws.Range("E4").Value = 0.5186869123087048
$ws.Range("G4").Value = 0.75
$ws.Range("I4").Value = 1

# Row 5
$ws.Range("E5").Value = 0.5582612362888708
$ws.Range("F5").Value = 0.7966863272106656
$ws.Range("G5").Value = 1

# Row 6
$ws.Range("D6").Value = "Only the village elder, Mira the Wise, was permitted to touch the stone with her bare hands."
$ws.Range("E6").Value = 0.4523436136211799
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 0.9999999999
$ws.Range("I6").Value = 1
